$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.313.81"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.865.82"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2859"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06570"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07836"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "1.865.62"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6989"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.095"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "269.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "30.276.13"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007634"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "2.110.96"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.231"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.941"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09923"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.354"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.049"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04728"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01873"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.766"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.311"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.949"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4172"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8343"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "970.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.116"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05680"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
